# Updates cryptos list figures (price/volume) for the Tue Jul 9 2024 run.
# Numeric-looking text values are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.573.81'
$ws.Range('E2').Value = '  +3.26%  '
$ws.Range('D3').Value = '3.081.15'
$ws.Range('E3').Value = '  +5.17%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''515.27'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').Value = '''142.73'
$ws.Range('E6').Value = '  +6.92%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '''0.436'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('D9').Value = '''7.28'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').Value = '  +3.90%  '
$ws.Range('E11').Value = '  +5.96%  '
$ws.Range('D12').Value = '3.605.32'
$ws.Range('E12').Value = '  +5.31%  '
$ws.Range('E13').Value = '  +2.86%  '
$ws.Range('D14').Value = '''25.79'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '''0.0000166'
$ws.Range('E15').Value = '  +3.00%  '
$ws.Range('D16').Value = '57.591.24'
$ws.Range('E16').Value = '  +3.43%  '
$ws.Range('D17').Value = '3.105.13'
$ws.Range('E17').Value = '  +6.11%  '
$ws.Range('D18').Value = '''6.06'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = '''13.07'
$ws.Range('E19').Value = '  +2.95%  '
$ws.Range('D20').Value = '''8.21'
$ws.Range('E20').Value = '  +6.40%  '
$ws.Range('D21').Value = '''339.59'
$ws.Range('E21').Value = '  +7.68%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = '''0.501'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = '''65.65'
$ws.Range('E24').Value = '  +4.47%  '
$ws.Range('E25').Value = '  +7.13%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0₃0947'
$ws.Range('E26').Value = '  +12.44%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '''6.48'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').Value = '''7.14'
$ws.Range('E29').Value = '  +4.01%  '
$ws.Range('D30').Value = '''1.82'
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('D31').Value = '''20.82'
$ws.Range('E31').Value = '  +4.63%  '
$ws.Range('D32').Value = '''1.19'
$ws.Range('E32').Value = '  +4.17%  '
$ws.Range('D33').Value = '''154.34'
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('E34').Value = '  +3.63%  '
$ws.Range('D35').Value = '''5.92'
$ws.Range('E35').Value = '  +5.27%  '
$ws.Range('D36').Value = '''26.33'
$ws.Range('E36').Value = '  +9.23%  '
$ws.Range('D37').Value = '''1.25'
$ws.Range('E37').Value = '  +4.38%  '
$ws.Range('D38').Value = '''0.0681'
$ws.Range('E38').Value = '  +5.37%  '
$ws.Range('D39').Value = '3.118.77'
$ws.Range('E39').Value = '  +5.38%  '
$ws.Range('D40').Value = '''37.12'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '''0.673'
$ws.Range('E41').Value = '  +5.61%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''3.87'
$ws.Range('E42').Value = '  +4.36%  '
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').Value = '2.269.23'
$ws.Range('E44').Value = '  +6.96%  '
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('D46').Value = '''0.0252'
$ws.Range('E46').Value = '  +7.23%  '
$ws.Range('D47').Value = '''0.960'
$ws.Range('E47').Value = '  +4.12%  '
$ws.Range('D48').Value = '''20.27'
$ws.Range('E48').Value = '  +8.17%  '
$ws.Range('D49').Value = '''5.88'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('D50').Value = '''0.0871'
$ws.Range('E50').Value = '  +4.00%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '''1.74'
$ws.Range('E51').Value = '  +3.58%  '
